$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column K
$ws.Range("K2").Value = "Ratio"

# ---- New rows of data (write column A first for all PUUM rows, then all
# TEAK02 rows, then go back and fill the "null" placeholder cells -- this
# mirrors the order the strings were first typed in and keeps the shared
# string table in the same order as the source file) ---------------------

$ws.Range("A9").Value = "PUUM"
$ws.Range("A10").Value = "PUUM"
$ws.Range("A11").Value = "PUUM"
$ws.Range("A12").Value = "TEAK02"
$ws.Range("A13").Value = "TEAK02"
$ws.Range("A14").Value = "TEAK02"

# Row 9: PUUM
$ws.Range("B9").Value = 43972
$ws.Range("B9").NumberFormat = "d-mmm-yy"
$ws.Range("C9").Value = 10
$ws.Range("D9").Value = 4.39
$ws.Range("E9").Value = 1.42
$ws.Range("F9").Value = 0.77
$ws.Range("G9").Value = 0.53
$ws.Range("H9").Value = 0.62
$ws.Range("I9").Value = 0.45
$ws.Range("J9").Formula = "=D9-E9"
$ws.Range("K9").Formula = "=J9/D9"

# Row 10: PUUM (nulls)
$ws.Range("B10").Value = 43972
$ws.Range("B10").NumberFormat = "d-mmm-yy"
$ws.Range("C10").Value = 39
$ws.Range("D10").Value = 1.14
$ws.Range("F10").Value = 0.3
$ws.Range("H10").Value = 0.175

# Row 11: PUUM (nulls)
$ws.Range("B11").Value = 43972
$ws.Range("B11").NumberFormat = "d-mmm-yy"
$ws.Range("C11").Value = 41
$ws.Range("D11").Value = 3.5
$ws.Range("F11").Value = 0.75
$ws.Range("H11").Value = 0.71

# Row 12: TEAK02
$ws.Range("B12").Value = 44341
$ws.Range("B12").NumberFormat = "d-mmm-yy"
$ws.Range("C12").Value = 47
$ws.Range("D12").Value = 3.11
$ws.Range("E12").Value = 1.03
$ws.Range("F12").Value = 0.6
$ws.Range("G12").Value = 0.47
$ws.Range("H12").Value = 0.26
$ws.Range("I12").Value = 0.41
$ws.Range("J12").Formula = "=D12-E12"
$ws.Range("K12").Formula = "=J12/D12"

# Row 13: TEAK02
$ws.Range("B13").Value = 44341
$ws.Range("B13").NumberFormat = "d-mmm-yy"
$ws.Range("C13").Value = 57
$ws.Range("D13").Value = 4.86
$ws.Range("E13").Value = 1.73
$ws.Range("F13").Value = 0.67
$ws.Range("G13").Value = 0.34
$ws.Range("H13").Value = 0.78
$ws.Range("I13").Value = 0.29
$ws.Range("J13").Formula = "=D13-E13"
$ws.Range("K13").Formula = "=J13/D13"

# Row 14: TEAK02
$ws.Range("B14").Value = 44341
$ws.Range("B14").NumberFormat = "d-mmm-yy"
$ws.Range("C14").Value = 46
$ws.Range("D14").Value = 1.55
$ws.Range("E14").Value = 1.22
$ws.Range("F14").Value = 0.42
$ws.Range("G14").Value = 0.35
$ws.Range("H14").Value = 0.26
$ws.Range("I14").Value = 0.24
$ws.Range("J14").Formula = "=D14-E14"
$ws.Range("K14").Formula = "=J14/D14"

# ---- Existing rows 3-8: add the new Ratio column (K) ---------------------
$ws.Range("K3").Formula = "=J3/D3"
$ws.Range("K4").Formula = "=J4/D4"
$ws.Range("K5").Formula = "=J5/D5"
$ws.Range("K6").Formula = "=J6/D6"
$ws.Range("K7").Formula = "=J7/D7"
$ws.Range("K8").Formula = "=J8/D8"

# Now go back and fill in the literal "null" placeholders for rows 10 & 11
$ws.Range("E10").Value = "null"
$ws.Range("G10").Value = "null"
$ws.Range("I10").Value = "null"
$ws.Range("J10").Value = "null"
$ws.Range("K10").Value = "null"

$ws.Range("E11").Value = "null"
$ws.Range("G11").Value = "null"
$ws.Range("I11").Value = "null"
$ws.Range("J11").Value = "null"
$ws.Range("K11").Value = "null"

# ---- Formatting -----------------------------------------------------------
$ws.Columns.Item(2).AutoFit() | Out-Null

# ---- Selection --------------------------------------------------------
$ws.Range("E15").Select()
